# Update gh-pages output data (generated at 456a3b4)
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8035
$ws1.Range("F10").Value = 478
$ws1.Range("F13").Value = 458
$ws1.Range("F17").Value = 5930
$ws1.Range("F18").Value = 190
$ws1.Range("F19").Value = 280
$ws1.Range("F20").Value = 1958
$ws1.Range("F21").Value = 35
$ws1.Range("F22").Value = 54

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 47

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8035
$ws4.Range("F10").Value = 478
$ws4.Range("F13").Value = 458
$ws4.Range("F18").Value = 5930
$ws4.Range("F19").Value = 47
$ws4.Range("F20").Value = 190
$ws4.Range("F21").Value = 280
$ws4.Range("F22").Value = 1958
$ws4.Range("F23").Value = 35
$ws4.Range("F24").Value = 54
